$wb = $excel.ActiveWorkbook

# ALC row 29: Dripping with Venom | Weak Blinding Potion
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 1466.5
$ws.Range("J29").Value = 2433
$ws.Range("L29").Value = 7299
$ws.Range("N29").Value = -7861

# ALC row 40: Stuck in the Moment | Horn Glue
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3573.963
$ws.Range("I40").Value = 2695
$ws.Range("K40").Value = 2695
$ws.Range("M40").Value = -2520

# ALC row 62: The Mustache Suits Him | Enchanted Mythrite Ink
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 2605.7856
$ws.Range("I62").Value = 2849.8
$ws.Range("J62").Value = 1995.75
$ws.Range("K62").Value = 2849.8
$ws.Range("L62").Value = 1995.75
$ws.Range("M62").Value = -2225.8
$ws.Range("N62").Value = -3243.75

# ALC row 65: Forgery of Convenience (L) | Enchanted Mythrite Ink
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 2605.7856
$ws.Range("I65").Value = 2849.8
$ws.Range("J65").Value = 1995.75
$ws.Range("K65").Value = 14249
$ws.Range("L65").Value = 9978.75
$ws.Range("M65").Value = -11129
$ws.Range("N65").Value = -16218.75

# ALC row 112: Making Ends Meet | Superior Spiritbond Potion
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 3485.9583
$ws.Range("I112").Value = 1566
$ws.Range("J112").Value = 3760.238
$ws.Range("K112").Value = 4698
$ws.Range("L112").Value = 11280.714
$ws.Range("M112").Value = -3590
$ws.Range("N112").Value = -13496.714

# ALC row 113: Amaro Kart | Starch Glue
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 5711.4443
$ws.Range("I113").Value = 4601
$ws.Range("J113").Value = 7099.5
$ws.Range("K113").Value = 4601
$ws.Range("L113").Value = 7099.5
$ws.Range("M113").Value = -1347
$ws.Range("N113").Value = -13607.5

# ALC row 116: Growing Up | Growth Formula Kappa
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 9648.666999999999
$ws.Range("I116").Value = 8581.333000000001
$ws.Range("K116").Value = 8581.333000000001
$ws.Range("M116").Value = -5139.333000000001

# ALC row 132: Fast-forwarding Flora | Growth Formula Lambda
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 2868.8572
$ws.Range("I132").Value = 2679.121
$ws.Range("K132").Value = 8037.363
$ws.Range("M132").Value = -5507.363

# ALC row 137: Cutting Edge of Culinary Quality | Magnesia Whetstone
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1799.1143
$ws.Range("I137").Value = 1470.8
$ws.Range("J137").Value = 2045.35
$ws.Range("K137").Value = 4412.4
$ws.Range("L137").Value = 6136.049999999999
$ws.Range("M137").Value = -1862.4
$ws.Range("N137").Value = -11236.05

# ALC row 138: All-night Crafting | Cunning Craftsman's Tisane
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2888.96
$ws.Range("J138").Value = 4018.5
$ws.Range("L138").Value = 12055.5
$ws.Range("N138").Value = -22335.5

# ARM row 61: Dealing with the Tough Stuff | Cobalt Ingot
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5071
$ws.Range("I61").Value = 4658.96
$ws.Range("J61").Value = 8504.666999999999
$ws.Range("K61").Value = 4658.96
$ws.Range("L61").Value = 8504.666999999999
$ws.Range("M61").Value = -4446.96
$ws.Range("N61").Value = -8928.666999999999

# ARM row 97: Ore for Me | High Steel Ingot
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 944.2857
$ws.Range("I97").Value = 277.25
$ws.Range("K97").Value = 277.25
$ws.Range("M97").Value = 218.75

# ARM row 102: Smells of Rich Tama-hagane | Tama-hagane Ingot
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 1134.8334
$ws.Range("I102").Value = 868.375
$ws.Range("K102").Value = 868.375
$ws.Range("M102").Value = 753.625

# ARM row 132: Don't Bore Me, Ore Me | Mountain Chromite Ingot
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 3674.5
$ws.Range("I132").Value = 2019.75
$ws.Range("K132").Value = 6059.25
$ws.Range("M132").Value = -3529.25

# ARM row 136: Metal with Mettle | Cobalt Tungsten Ingot
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 5071
$ws.Range("I136").Value = 4658.96
$ws.Range("J136").Value = 8504.666999999999
$ws.Range("K136").Value = 13976.88
$ws.Range("L136").Value = 25514.001
$ws.Range("M136").Value = -11426.88
$ws.Range("N136").Value = -30614.001

# BSM row 20: Smelt and Dealt | Iron Ingot
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1858.909
$ws.Range("I20").Value = 2940
$ws.Range("J20").Value = 1241.1428
$ws.Range("K20").Value = 2940
$ws.Range("L20").Value = 1241.1428
$ws.Range("M20").Value = -2693
$ws.Range("N20").Value = -1735.1428

# BSM row 86: Through Thick and Thin | Adamantite Nugget
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2572.625
$ws.Range("I86").Value = 2654.4285
$ws.Range("J86").Value = 2000
$ws.Range("K86").Value = 2654.4285
$ws.Range("L86").Value = 2000
$ws.Range("M86").Value = -1531.4285
$ws.Range("N86").Value = -4246

# BSM row 89: Piercing Eyes Deserve Piercing Shafts (L) | Adamantite Nugget
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 2572.625
$ws.Range("I89").Value = 2654.4285
$ws.Range("J89").Value = 2000
$ws.Range("K89").Value = 13272.1425
$ws.Range("L89").Value = 10000
$ws.Range("M89").Value = -7656.1425
$ws.Range("N89").Value = -21232

# BSM row 107: The Gold Experience | Deepgold Nugget
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 4012.9412
$ws.Range("I107").Value = 4236
$ws.Range("J107").Value = 3604
$ws.Range("K107").Value = 4236
$ws.Range("L107").Value = 3604
$ws.Range("M107").Value = -2316
$ws.Range("N107").Value = -7444

# CRP row 7: Gridania's Got Talent | Maple Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 258.5
$ws.Range("I7").Value = 169.22223
$ws.Range("J7").Value = 347.77777
$ws.Range("K7").Value = 169.22223
$ws.Range("L7").Value = 347.77777
$ws.Range("M7").Value = -56.22223
$ws.Range("N7").Value = -573.7777699999999

# CRP row 31: Wall Not Found | Walnut Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7758.364
$ws.Range("I31").Value = 3693.2942
$ws.Range("J31").Value = 21579.6
$ws.Range("K31").Value = 3693.2942
$ws.Range("L31").Value = 21579.6
$ws.Range("M31").Value = -3398.2942
$ws.Range("N31").Value = -22169.6

# CRP row 34: Armoires of the Rich and Famous | Walnut Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 7758.364
$ws.Range("I34").Value = 3693.2942
$ws.Range("J34").Value = 21579.6
$ws.Range("K34").Value = 3693.2942
$ws.Range("L34").Value = 21579.6
$ws.Range("M34").Value = -3491.2942
$ws.Range("N34").Value = -21983.6

# CRP row 58: You Do the Heavy Lifting | Mahogany Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3940.6843
$ws.Range("I58").Value = 1389
$ws.Range("J58").Value = 5429.1665
$ws.Range("K58").Value = 1389
$ws.Range("L58").Value = 5429.1665
$ws.Range("M58").Value = -1186
$ws.Range("N58").Value = -5835.1665

# CRP row 134: Wood You Be Quiet | Ceiba Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 3936.4358
$ws.Range("J134").Value = 8398.799999999999
$ws.Range("L134").Value = 25196.4
$ws.Range("N134").Value = -30266.4

# CRP row 136: Turali Quality | Dark Mahogany Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 3940.6843
$ws.Range("I136").Value = 1389
$ws.Range("J136").Value = 5429.1665
$ws.Range("K136").Value = 4167
$ws.Range("L136").Value = 16287.4995
$ws.Range("M136").Value = -1617
$ws.Range("N136").Value = -21387.4995

# CRP row 139: Weaving a Path | Acacia Spinning Wheel
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H139").Value = 119996.5
$ws.Range("J139").Value = 119996.5
$ws.Range("L139").Value = 119996.5
$ws.Range("N139").Value = -130276.5

# CUL row 2: Pork Is a Salty Food | Table Salt
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 522
$ws.Range("I2").Value = 35.6
$ws.Range("K2").Value = 213.6
$ws.Range("M2").Value = -100.6

# CUL row 113: Can't Eat Just One | Night Vinegar
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 3321
$ws.Range("J113").Value = 2611.3333
$ws.Range("L113").Value = 7833.999899999999
$ws.Range("N113").Value = -12173.9999

# CUL row 118: Teetotally | Masala Chai
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H118").Value = 6366
$ws.Range("I118").Value = 6366
$ws.Range("K118").Value = 19098
$ws.Range("M118").Value = -17855

# CUL row 132: More Mezcal | Cooking Mezcal
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 2633.5293
$ws.Range("I132").Value = 2588.3635
$ws.Range("J132").Value = 2716.3333
$ws.Range("K132").Value = 23295.2715
$ws.Range("L132").Value = 24446.9997
$ws.Range("M132").Value = -20765.2715
$ws.Range("N132").Value = -29506.9997

# CUL row 140: Sweet, Sweet Bean Juice | Mesquite Juice
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 1609.9
$ws.Range("I140").Value = 1588.5555
$ws.Range("J140").Value = 1802
$ws.Range("K140").Value = 4765.666499999999
$ws.Range("L140").Value = 5406
$ws.Range("M140").Value = 414.3335000000006
$ws.Range("N140").Value = -15766

# GSM row 31: One and Only | Staghorn Staff
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H31").Value = 7000
$ws.Range("I31").Value = 7666.6665
$ws.Range("K31").Value = 7666.6665
$ws.Range("M31").Value = -7374.6665

# GSM row 37: Dancing with the Stars | Toothed Staghorn Staff
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H37").Value = 7000
$ws.Range("I37").Value = 7666.6665
$ws.Range("K37").Value = 7666.6665
$ws.Range("M37").Value = -7389.6665

# GSM row 70: Sky Is the Limit | Mythrite Ingot
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5598.919
$ws.Range("I70").Value = 5347.5
$ws.Range("J70").Value = 5613.2856
$ws.Range("K70").Value = 5347.5
$ws.Range("L70").Value = 5613.2856
$ws.Range("M70").Value = -5077.5
$ws.Range("N70").Value = -6153.2856

# GSM row 73: Hulls of Broken Dreams (L) | Mythrite Ingot
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 5598.919
$ws.Range("I73").Value = 5347.5
$ws.Range("J73").Value = 5613.2856
$ws.Range("K73").Value = 5347.5
$ws.Range("L73").Value = 5613.2856
$ws.Range("M73").Value = -4411.5
$ws.Range("N73").Value = -7485.2856

# GSM row 132: On Board for Lar | Lar Ingot
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4359.0356
$ws.Range("I132").Value = 3003.261
$ws.Range("J132").Value = 10595.6
$ws.Range("K132").Value = 9009.782999999999
$ws.Range("L132").Value = 31786.8
$ws.Range("M132").Value = -6479.782999999999
$ws.Range("N132").Value = -36846.8

# LTW row 7: Tan Before the Ban | Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 44820.383
$ws.Range("I7").Value = 44820.383
$ws.Range("K7").Value = 44820.383
$ws.Range("M7").Value = -44708.383

# LTW row 55: It's Not a Job, It's a Calling | Peiste Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 1062.375
$ws.Range("I55").Value = 269.42856
$ws.Range("K55").Value = 269.42856
$ws.Range("M55").Value = -96.42856

# LTW row 61: Spelling Me Softly | Raptor Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 5584.8237
$ws.Range("I61").Value = 2869.2222
$ws.Range("K61").Value = 2869.2222
$ws.Range("M61").Value = -2667.2222

# LTW row 82: Trainin' the Neck | Dragon Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1566.8
$ws.Range("I82").Value = 502
$ws.Range("K82").Value = 502
$ws.Range("M82").Value = -141

# LTW row 85: Training Is Only Skintight (L) | Dragon Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 1566.8
$ws.Range("I85").Value = 502
$ws.Range("K85").Value = 502
$ws.Range("M85").Value = 746

# LTW row 113: Peace in Rest | Atrociraptor Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 5584.8237
$ws.Range("I113").Value = 2869.2222
$ws.Range("K113").Value = 2869.2222
$ws.Range("M113").Value = -699.2222000000002

# LTW row 126: Battered Books | Saiga Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 44820.383
$ws.Range("I126").Value = 44820.383
$ws.Range("K126").Value = 134461.149
$ws.Range("M126").Value = -131991.149

# LTW row 132: Tenets of Tanning | Silver Lobo Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 6127.1816
$ws.Range("I132").Value = 3387.5
$ws.Range("J132").Value = 13433
$ws.Range("K132").Value = 10162.5
$ws.Range("L132").Value = 40299
$ws.Range("M132").Value = -7632.5
$ws.Range("N132").Value = -45359

# WVR row 97: Getting a Leg Up | Ruby Cotton Gaskins of Striking
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H97").Value = 76369.75
$ws.Range("J97").Value = 76369.75
$ws.Range("L97").Value = 76369.75
$ws.Range("N97").Value = -78351.75

# WVR row 113: A Tender Table | Pixie Floss
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1112.4
$ws.Range("J113").Value = 1356.6666
$ws.Range("L113").Value = 4069.9998
$ws.Range("N113").Value = -8409.9998

# WVR row 132: Comfy Cabins | Snow Cotton Cloth
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3504
$ws.Range("I132").Value = 2368.7273
$ws.Range("K132").Value = 7106.1819
$ws.Range("M132").Value = -4576.1819
